# Update "想去人数" (interested-count) figures that changed between the two
# data pulls recorded in this workbook. The same five events appear on the
# "展览" sheet (rows 2-6) and on the "全部类型" sheet (rows 4-8); update both.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1040
$ws1.Range("F3").Value = 196
$ws1.Range("F4").Value = 2406
$ws1.Range("F5").Value = 29
$ws1.Range("F6").Value = 528

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1040
$ws4.Range("F5").Value = 196
$ws4.Range("F6").Value = 2406
$ws4.Range("F7").Value = 29
$ws4.Range("F8").Value = 528
